$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FilePath column (B) for each scene row: these values used to be
# blank (rows 3-7) or point at an unused/stale path (row 2), and now point at
# the per-scene ini files.
$ws.Range("B2").Value = "../../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("B3").Value = "../../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("B4").Value = "../../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("B5").Value = "../../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("B6").Value = "../../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("B7").Value = "../../NFDataCfg/Ini/Scene/6.xml"

# The newly-typed values (rows 3-7) picked up the text number format plus a
# distinct "family 3" flavour of the 宋体 font, same as when typed directly
# in the Excel UI.
$ws.Range("B3:B7").NumberFormat = "@"
$ws.Range("B3:B7").Font.Family = 3

# Selection moved off E2 (and the frozen/scrolled-to topLeftCell) onto B5.
$ws.Range("B5").Select() | Out-Null
